$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the room ("Salle") column F for each course row with the
# corresponding room name (previously empty cells).
$ws.Range("F3").Value  = "U3-110"
$ws.Range("F6").Value  = "U3-110"
$ws.Range("F8").Value  = "U3-110"
$ws.Range("F11").Value = "U3-109"
$ws.Range("F13").Value = "U3-107"
$ws.Range("F16").Value = "U3-110"
$ws.Range("F18").Value = "U3-110"
$ws.Range("F21").Value = "U3-Amphi"
